# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Most "Price" cells are plain numeric-looking text (e.g. "147.50") that the
# Excel object model would otherwise auto-convert to a real number on
# assignment; a leading apostrophe forces those to stay literal text, just
# like the original inline strings in the sheet. Cells whose text already
# can't parse as a single number (thousands-separated "67.438.27", the
# subscript-digit PEPE price, etc.) are set without the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.438.27"
$ws.Range("E2").Value = "  -2.44%  "
$ws.Range("D3").Value = "3.261.92"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'587.25"
$ws.Range("E5").Value = "  -3.90%  "
$ws.Range("D6").Value = "'147.50"
$ws.Range("E6").Value = "  -11.41%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.257.05"
$ws.Range("E8").Value = "  -4.89%  "
$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  -8.97%  "
$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  -12.54%  "
$ws.Range("D11").Value = "'6.60"
$ws.Range("E11").Value = "  -5.08%  "
$ws.Range("D12").Value = "'0.499"
$ws.Range("E12").Value = "  -11.06%  "
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "  -9.08%  "
$ws.Range("D14").Value = "'37.71"
$ws.Range("E14").Value = "  -14.25%  "
$ws.Range("D15").Value = "3.787.74"
$ws.Range("E15").Value = "  -5.22%  "
$ws.Range("D16").Value = "67.540.86"
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "3.272.49"
$ws.Range("E17").Value = "  -5.01%  "
$ws.Range("D18").Value = "'0.113"
$ws.Range("E18").Value = "  -6.13%  "
$ws.Range("D19").Value = "'521.66"
$ws.Range("E19").Value = "  -9.96%  "
$ws.Range("D20").Value = "'7.02"
$ws.Range("E20").Value = "  -13.17%  "
$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  -13.46%  "
$ws.Range("D22").Value = "'0.745"
$ws.Range("E22").Value = "  -11.35%  "
$ws.Range("D23").Value = "'7.71"
$ws.Range("E23").Value = "  -12.48%  "
$ws.Range("D24").Value = "'84.92"
$ws.Range("E24").Value = "  -11.01%  "
$ws.Range("D25").Value = "'13.27"
$ws.Range("E25").Value = "  -11.69%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'3.19"
$ws.Range("E27").Value = "  -11.57%  "
$ws.Range("D28").Value = "'2.12"
$ws.Range("E28").Value = "  -11.49%  "
$ws.Range("D29").Value = "'7.89"
$ws.Range("E29").Value = "  -7.34%  "
$ws.Range("D30").Value = "'28.66"
$ws.Range("E30").Value = "  -11.98%  "
$ws.Range("D31").Value = "'1.18"
$ws.Range("E31").Value = "  -4.65%  "
$ws.Range("D32").Value = "'2.62"
$ws.Range("E32").Value = "  -4.81%  "
$ws.Range("D33").Value = "'6.47"
$ws.Range("E33").Value = "  -16.58%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'57.14"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'5.61"
$ws.Range("E35").Value = "  -13.99%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'504.01"
$ws.Range("E37").Value = "  -12.39%  "
$ws.Range("D38").Value = "'0.0440"
$ws.Range("E38").Value = "  -6.10%  "
$ws.Range("D39").Value = "'0.0841"
$ws.Range("E39").Value = "  -11.41%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.124"
$ws.Range("E40").Value = "  -11.87%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").Value = "'8.77"
$ws.Range("E41").Value = "  -15.66%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.919.11"
$ws.Range("E42").Value = "  -9.45%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.68"
$ws.Range("E43").Value = "  -13.48%  "
$ws.Range("D44").Value = "'0.263"
$ws.Range("E44").Value = "  -10.10%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.17"
$ws.Range("E45").Value = "  -8.31%  "
$ws.Range("B46").Value = "PEPE"
$ws.Range("C46").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D46").Value = "0.0₃0564"
$ws.Range("E46").Value = "  -16.52%  "
$ws.Range("D47").Value = "'26.23"
$ws.Range("E47").Value = "  -15.28%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'123.86"
$ws.Range("E49").Value = "  -6.59%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  -17.10%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.112"
$ws.Range("E51").Value = "  -10.50%  "
